# Apply updated cryptocurrency price/volume data to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.789.49"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.447.67"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.438.98"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "2.891.10"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "61.753.90"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "2.441.90"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "583.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.28%  "
$ws.Range("D29").Value = "2.565.82"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "0.0₃0925"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  -4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  +26.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -2.02%  "
